$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Formula = "'1"
$ws.Range("E3").Formula = "'314"
$ws.Range("F3").Formula = "'4"
$ws.Range("G3").Formula = "'3"
$ws.Range("E7").Formula = "'1147"
$ws.Range("F7").Formula = "'13"
$ws.Range("G7").Formula = "'13"
$ws.Range("I7").Formula = "'1"
$ws.Range("L7").Formula = "'7"
$ws.Range("E8").Formula = "'1013"
$ws.Range("F8").Formula = "'14"
$ws.Range("G8").Formula = "'10"
$ws.Range("E9").Formula = "'377"
$ws.Range("F9").Formula = "'7"
$ws.Range("H9").Formula = "'3"
$ws.Range("J9").Formula = "'9"
$ws.Range("E10").Formula = "'687"
$ws.Range("F10").Formula = "'10"
$ws.Range("G10").Formula = "'8"
$ws.Range("L10").Formula = "'3"
$ws.Range("E11").Formula = "'229"
$ws.Range("F11").Formula = "'5"
$ws.Range("H11").Formula = "'3"
$ws.Range("J11").Formula = "'6"
$ws.Range("E12").Formula = "'1239"
$ws.Range("F12").Formula = "'14"
$ws.Range("G12").Formula = "'14"
$ws.Range("E14").Formula = "'181"
$ws.Range("F14").Formula = "'4"
$ws.Range("G14").Formula = "'2"
$ws.Range("I14").Formula = "'1"
$ws.Range("J15").Formula = "'10"
$ws.Range("E16").Formula = "'707"
$ws.Range("F16").Formula = "'14"
$ws.Range("H16").Formula = "'5"
$ws.Range("J16").Formula = "'5"
$ws.Range("E17").Formula = "'396"
$ws.Range("F17").Formula = "'12"
$ws.Range("H17").Formula = "'10"
$ws.Range("J17").Formula = "'10"
$ws.Range("E18").Formula = "'854"
$ws.Range("F18").Formula = "'13"
$ws.Range("G18").Formula = "'11"
$ws.Range("I18").Formula = "'9"
$ws.Range("J20").Formula = "'11"
$ws.Range("E23").Formula = "'278"
$ws.Range("F23").Formula = "'9"
$ws.Range("H23").Formula = "'6"
$ws.Range("J23").Formula = "'11"
$ws.Range("E25").Formula = "'1051"
$ws.Range("F25").Formula = "'13"
$ws.Range("G25").Formula = "'13"
$ws.Range("I25").Formula = "'7"
$ws.Range("E26").Formula = "'1081"
$ws.Range("F26").Formula = "'14"
$ws.Range("G26").Formula = "'12"
$ws.Range("E27").Formula = "'831"
$ws.Range("F27").Formula = "'14"
$ws.Range("G27").Formula = "'10"
$ws.Range("E29").Formula = "'407"
$ws.Range("F29").Formula = "'12"
$ws.Range("G29").Formula = "'3"
$ws.Range("I29").Formula = "'3"
$ws.Range("J30").Formula = "'1"
